$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop every existing hyperlink up front. (Item-level .Delete() is a
#    no-op in this engine, but calling .Delete() on a range-scoped
#    Hyperlinks collection clears the whole worksheet collection - so we
#    use that, then rebuild every hyperlink explicitly once the data is
#    in its final place.)
# ------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Make room for the new "Q_IAB" source row at row 37 by shifting the
#    existing rows 37..56 down to 38..57. Only cell VALUES are copied
#    (bottom-up, so nothing is clobbered before it's read) which keeps
#    every destination cell's pre-existing style (every data row already
#    shares the same per-column style, so nothing needs reformatting).
# ------------------------------------------------------------------
for ($r = 56; $r -ge 37; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 12; $c++) {
        $val = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($dest, $c).Value = $val
    }
}

# Row 57 is brand new territory (the sheet used to stop at 56), so it has
# no style of its own yet - clone the formatting from row 56 (which now
# also holds correctly-shifted data) so the new last row matches every
# other data row instead of picking up a freshly-minted style.
$ws.Range("A56:L56").Copy()
$ws.Range("A57:L57").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3. Write the new Q_IAB row into the now-empty row 37.
# ------------------------------------------------------------------
$ws.Cells.Item(37,1).Value  = "Q_IAB"
$ws.Cells.Item(37,2).Value  = "Institut für Arbeitsmarkt- und Berufsforschung (IAB)"
$ws.Cells.Item(37,3).Value  = "des Instituts für Arbeitsmarkt- und Berufsforschung"
$ws.Cells.Item(37,4).Value  = "Institut for Employment Research"
$ws.Cells.Item(37,5).Value  = "the Institut for Employment Research"
$ws.Cells.Item(37,6).Value  = "Institut für Arbeitsmarkt- und Berufsforschung"
$ws.Cells.Item(37,7).Value  = "Institut for Employment Research"
$ws.Cells.Item(37,8).Value  = "https://iab.de/"
$ws.Cells.Item(37,9).Value  = "https://iab.de/en/"
$ws.Cells.Item(37,10).Value = ""
$ws.Cells.Item(37,11).Value = ""
$ws.Cells.Item(37,12).Value = "iab"

# ------------------------------------------------------------------
# 4. Re-create all the hyperlinks. Rows above the insertion point (< 37)
#    keep their original addresses; the Stifterverband row (previously
#    row 48) now lives at row 49; and the new IAB row gets its own pair
#    of links at row 37.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("H3"),  "https://ag-energiebilanzen.de/")
$ws.Hyperlinks.Add($ws.Range("I3"),  "https://ag-energiebilanzen.de/en/")
$ws.Hyperlinks.Add($ws.Range("I8"),  "https://www.bbsr.bund.de/BBSR/EN/home/_node.html")
$ws.Hyperlinks.Add($ws.Range("H12"), "https://www.kulturstaatsministerin.de/DE/startseite/startseite_node.html")
$ws.Hyperlinks.Add($ws.Range("I12"), "https://www.kulturstaatsministerin.de/DE/startseite/startseite_node.html")
$ws.Hyperlinks.Add($ws.Range("H14"), "https://bmdv.bund.de/DE/Home/home.html")
$ws.Hyperlinks.Add($ws.Range("I14"), "https://bmdv.bund.de/EN/Home/home.html")
$ws.Hyperlinks.Add($ws.Range("H17"), "https://www.bmfsfj.de/")
$ws.Hyperlinks.Add($ws.Range("I17"), "https://www.bmfsfj.de/en")
$ws.Hyperlinks.Add($ws.Range("H20"), "https://www.bmz.de/de")
$ws.Hyperlinks.Add($ws.Range("I20"), "https://www.bmz.de/en")
$ws.Hyperlinks.Add($ws.Range("H24"), "https://ceval.de/")
$ws.Hyperlinks.Add($ws.Range("I24"), "https://ceval.de/en/")
$ws.Hyperlinks.Add($ws.Range("H25"), "https://www.deutsche-digitale-bibliothek.de/")
$ws.Hyperlinks.Add($ws.Range("I25"), "https://www.deutsche-digitale-bibliothek.de/?lang=en")
$ws.Hyperlinks.Add($ws.Range("H35"), "https://ghsindex.org/")
$ws.Hyperlinks.Add($ws.Range("I35"), "https://ghsindex.org/")
$ws.Hyperlinks.Add($ws.Range("H37"), "https://iab.de/")
$ws.Hyperlinks.Add($ws.Range("I37"), "https://iab.de/en/")
$ws.Hyperlinks.Add($ws.Range("H49"), "https://www.stifterverband.org/")
$ws.Hyperlinks.Add($ws.Range("I49"), "https://www.stifterverband.org/english")
